# Update the "example" worksheet with the refreshed run results.
$wb = $excel.ActiveWorkbook
$wsExample = $wb.Worksheets.Item("example")
$wsIndicators = $wb.Worksheets.Item("indicators")

# --- "example" sheet, row 2 (totals row) ---
$wsExample.Range("B2").Value = 128535.12
$wsExample.Range("C2").Value = 89974.58
$wsExample.Range("D2").Value = 38560.54
$wsExample.Range("E2").Value = 85905.53
$wsExample.Range("F2").Value = 5694.52
$wsExample.Range("G2").Value = 68505.96000000001
$wsExample.Range("H2").Value = 401.57
$wsExample.Range("I2").Value = 95.95999999999999
$wsExample.Range("J2").Value = 74368.59
$wsExample.Range("K2").Value = 73139.78999999999
$wsExample.Range("L2").Value = 73414.06
$wsExample.Range("M2").Value = 148400.96
$wsExample.Range("N2").Value = 21918.92
$wsExample.Range("O2").Value = 6025.53

# --- "example" sheet, row 4 (ratios) ---
$wsExample.Range("E4").Value = 1.012667396236561
$wsExample.Range("F4").Value = 1.157060676825656
$wsExample.Range("J4").Value = 1.002101752710336
$wsExample.Range("K4").Value = 1.000393627332655
$wsExample.Range("L4").Value = 1.034459229562604

# --- "example" sheet, row 8 ---
$wsExample.Range("K8").Value = 4.409754900570112
$wsExample.Range("L8").Value = 17.98259849983794

# --- "example" sheet, row 9 ---
$wsExample.Range("D9").Value = 36.33336666666666
$wsExample.Range("J9").Value = 6.687943893034388
$wsExample.Range("L9").Value = 25.35213050123805

# --- "example" sheet, row 10 ---
$wsExample.Range("E10").Value = 0.6049708395891255
$wsExample.Range("F10").Value = 0.6138556679281634
$wsExample.Range("G10").Value = 0.3534032273528521
$wsExample.Range("L10").Value = 0.3773914661224246

# --- "example" sheet, row 11 ---
$wsExample.Range("E11").Value = 0.04245409400625442
$wsExample.Range("F11").Value = 0.04307759073180095
$wsExample.Range("L11").Value = 0.003341406509317179

# --- "example" sheet, row 12 ---
$wsExample.Range("E12").Value = 0.03184057050469082
$wsExample.Range("F12").Value = 0.03230819304885071
$wsExample.Range("L12").Value = 0.03583579306542615

# --- "example" sheet, row 13 ---
$wsExample.Range("D13").Value = 10.78026666666666
$wsExample.Range("E13").Value = 0.1061352350156361
$wsExample.Range("F13").Value = 0.1076939768295024
$wsExample.Range("L13").Value = 5.85565126756877

# --- "example" sheet, single summary cells ---
$wsExample.Range("B16").Value = 74.96931608889483
$wsExample.Range("B17").Value = 15.0270094682561
$wsExample.Range("B18").Value = 1448.738524629777
$wsExample.Range("B19").Value = 703.3036818639184
$wsExample.Range("B26").Value = 401.5658607552735
$wsExample.Range("B27").Value = 95.96174510427927
$wsExample.Range("B28").Value = 74670.96155344439
$wsExample.Range("B29").Value = 74368.59279053439

# --- "indicators" sheet ---
$wsIndicators.Range("B3").Value = -148400.962485911
$wsIndicators.Range("B5").Value = 2242.038532050846
$wsIndicators.Range("B7").Value = 4439236.293460676
$wsIndicators.Range("B8").Value = 148400.962485911
$wsIndicators.Range("B9").Value = 3714921.257162649
$wsIndicators.Range("B10").Value = 6017025.973444377
